$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (Tipo), shifting it to E
$ws.Columns.Item(4).Insert()

$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D1").Value = "MAE"
$ws.Range("D2").Value = 0.2072673588334912

Write-Output "done"
